$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (2-11) with the new measurement set and
#     extend the table down to row 41 with the additional dated readings ---
$ws.Cells.Item(2,1).Value = 44564
$ws.Cells.Item(2,2).Value = 24.5
$ws.Cells.Item(2,3).Value = 25.01
$ws.Cells.Item(3,1).Value = 44565
$ws.Cells.Item(3,2).Value = 25.49
$ws.Cells.Item(3,3).Value = 25.5
$ws.Cells.Item(4,1).Value = 44566
$ws.Cells.Item(4,2).Value = 29.12
$ws.Cells.Item(4,3).Value = 28.73
$ws.Cells.Item(5,1).Value = 44567
$ws.Cells.Item(5,2).Value = 27.16
$ws.Cells.Item(5,3).Value = 28.24
$ws.Cells.Item(6,1).Value = 44568
$ws.Cells.Item(6,2).Value = 25
$ws.Cells.Item(6,3).Value = 24.69
$ws.Cells.Item(7,1).Value = 44569
$ws.Cells.Item(7,2).Value = 25
$ws.Cells.Item(7,3).Value = 26.08
$ws.Cells.Item(8,1).Value = 44571
$ws.Cells.Item(8,2).Value = 26.96
$ws.Cells.Item(8,3).Value = 26.86
$ws.Cells.Item(9,1).Value = 44572
$ws.Cells.Item(9,2).Value = 29.22
$ws.Cells.Item(9,3).Value = 28.73
$ws.Cells.Item(10,1).Value = 44573
$ws.Cells.Item(10,2).Value = 23.14
$ws.Cells.Item(10,3).Value = 22.11
$ws.Cells.Item(11,1).Value = 44574
$ws.Cells.Item(11,2).Value = 23.73
$ws.Cells.Item(11,3).Value = 23.82
$ws.Cells.Item(12,1).Value = 44575
$ws.Cells.Item(12,2).Value = 21.67
$ws.Cells.Item(12,3).Value = 21.76
$ws.Cells.Item(13,1).Value = 44590
$ws.Cells.Item(13,2).Value = 27.84
$ws.Cells.Item(13,3).Value = 27.25
$ws.Cells.Item(14,1).Value = 44591
$ws.Cells.Item(14,2).Value = 26.5
$ws.Cells.Item(14,3).Value = 25.8
$ws.Cells.Item(15,1).Value = 44592
$ws.Cells.Item(15,2).Value = 28.73
$ws.Cells.Item(15,3).Value = 28.53
$ws.Cells.Item(16,1).Value = 44593
$ws.Cells.Item(16,2).Value = 22.4
$ws.Cells.Item(16,3).Value = 20.8
$ws.Cells.Item(17,1).Value = 44594
$ws.Cells.Item(17,2).Value = 26.18
$ws.Cells.Item(17,3).Value = 26.57
$ws.Cells.Item(18,1).Value = 44595
$ws.Cells.Item(18,2).Value = 25.7
$ws.Cells.Item(18,3).Value = 24.8
$ws.Cells.Item(19,1).Value = 44596
$ws.Cells.Item(19,2).Value = 26.37
$ws.Cells.Item(19,3).Value = 26.86
$ws.Cells.Item(20,1).Value = 44601
$ws.Cells.Item(20,2).Value = 26.47
$ws.Cells.Item(20,3).Value = 26.67
$ws.Cells.Item(21,1).Value = 44602
$ws.Cells.Item(21,2).Value = 27.64
$ws.Cells.Item(21,3).Value = 28.04
$ws.Cells.Item(22,1).Value = 44603
$ws.Cells.Item(22,2).Value = 28.55
$ws.Cells.Item(22,3).Value = 29.1
$ws.Cells.Item(23,1).Value = 44604
$ws.Cells.Item(23,2).Value = 29.56
$ws.Cells.Item(23,3).Value = 20.6
$ws.Cells.Item(24,1).Value = 44606
$ws.Cells.Item(24,2).Value = 28.14
$ws.Cells.Item(24,3).Value = 28.92
$ws.Cells.Item(25,1).Value = 44607
$ws.Cells.Item(25,2).Value = 27.83
$ws.Cells.Item(25,3).Value = 28.84
$ws.Cells.Item(26,1).Value = 44609
$ws.Cells.Item(26,2).Value = 22.8
$ws.Cells.Item(26,3).Value = 21.9
$ws.Cells.Item(27,1).Value = 44610
$ws.Cells.Item(27,2).Value = 23.86
$ws.Cells.Item(27,3).Value = 24
$ws.Cells.Item(28,1).Value = 44611
$ws.Cells.Item(28,2).Value = 23.22
$ws.Cells.Item(28,3).Value = 23.39
$ws.Cells.Item(29,1).Value = 44613
$ws.Cells.Item(29,2).Value = 24.9
$ws.Cells.Item(29,3).Value = 24.31
$ws.Cells.Item(30,1).Value = 44614
$ws.Cells.Item(30,2).Value = 31.27
$ws.Cells.Item(30,3).Value = 31.47
$ws.Cells.Item(31,1).Value = 44615
$ws.Cells.Item(31,2).Value = 27.16
$ws.Cells.Item(31,3).Value = 28.57
$ws.Cells.Item(32,1).Value = 44616
$ws.Cells.Item(32,2).Value = 29.05
$ws.Cells.Item(32,3).Value = 29.8
$ws.Cells.Item(33,1).Value = 44617
$ws.Cells.Item(33,2).Value = 20.38
$ws.Cells.Item(33,3).Value = 21.5
$ws.Cells.Item(34,1).Value = 44618
$ws.Cells.Item(34,2).Value = 27.02
$ws.Cells.Item(34,3).Value = 28.2
$ws.Cells.Item(35,1).Value = 44621
$ws.Cells.Item(35,2).Value = 28.04
$ws.Cells.Item(35,3).Value = 27.9
$ws.Cells.Item(36,1).Value = 44622
$ws.Cells.Item(36,2).Value = 26.47
$ws.Cells.Item(36,3).Value = 27.5
$ws.Cells.Item(37,1).Value = 44623
$ws.Cells.Item(37,2).Value = 27.18
$ws.Cells.Item(37,3).Value = 27.08
$ws.Cells.Item(38,1).Value = 44624
$ws.Cells.Item(38,2).Value = 22.04
$ws.Cells.Item(38,3).Value = 22.53
$ws.Cells.Item(39,1).Value = 44625
$ws.Cells.Item(39,2).Value = 25.88
$ws.Cells.Item(39,3).Value = 25.88
$ws.Cells.Item(40,1).Value = 44626
$ws.Cells.Item(40,2).Value = 32.08
$ws.Cells.Item(40,3).Value = 31.37
$ws.Cells.Item(41,1).Value = 44627
$ws.Cells.Item(41,2).Value = 24.5
$ws.Cells.Item(41,3).Value = 25.3

# --- Make sure every date cell in column A (old + new rows) keeps the
#     original short-date number format by copying A2's style over the range ---
$ws.Cells.Item(2,1).Copy()
$ws.Range("A2:A41").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Move the selection/view down to the new last row, like the saved workbook ---
$ws.Range("A41").Select()

# --- Page setup was touched (portrait orientation) when the sheet was edited ---
$ws.PageSetup.Orientation = 1

